# "Fix import sur les utilisateurs"
#
# The "Site Web" column (I) on the "Liste des utilisateurs" sheet held bare
# host names ("www.google.fr") instead of full URLs, which broke the
# hyperlink import. Add the missing "http://" scheme, then leave the sheet
# and selection the way the user left them after making the fix.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Liste des utilisateurs")

# Make this the active sheet (it becomes the workbook's active tab).
$ws.Activate() | Out-Null

# Column I ("Site Web"), rows 4-19, all contained the literal text
# "www.google.fr" - rewrite them with an explicit http:// scheme.
$ws.Range("I4:I19").Value = "http://www.google.fr"

# Leave the selection where the user ended up: I9:I19, active cell I9.
$ws.Range("I9:I19").Select() | Out-Null
